# Applies the data corrections described by the diff:
#  - Fix mojibake in the footnote text (A103): literal ">" characters that
#    should have been accented Portuguese/Spanish letters (í / ú).
#  - Update a handful of recalculated numeric cells (J67, Z70, G73, J73, Z73,
#    and the full "Africa, Fragile States" / "ROW, Fragile States" aggregate
#    rows 97-98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footnote text correction (shared string used only by A103) ---------
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# --- Scattered floating-point refinements --------------------------------
$ws.Range("J67").Value = 480.13494270018799
$ws.Range("Z70").Value = 21.416995666564102
$ws.Range("G73").Value = 3.3716196462049002
$ws.Range("J73").Value = 532.448906333688
$ws.Range("Z73").Value = 33.903273713676597

# --- Row 97: "Africa, Fragile States" aggregate recalculation -----------
$ws.Range("C97").Value = 151565
$ws.Range("D97").Value = 194.53015348871301
$ws.Range("E97").Value = 4.9890148902285496
$ws.Range("F97").Value = 99.856380264071007
$ws.Range("G97").Value = 2.8724759546599001
$ws.Range("H97").Value = 6.3505367373854504
$ws.Range("I97").Value = 502570
$ws.Range("J97").Value = 645.036909832893
$ws.Range("K97").Value = 16.920450870827398
$ws.Range("L97").Value = 238.89370811806299
$ws.Range("M97").Value = 6.7360112134466297
$ws.Range("N97").Value = 14.857275743775601
$ws.Range("O97").Value = 6.2813352104625801
$ws.Range("P97").Value = 52.924163752232197
$ws.Range("Q97").Value = 28.780682623500802
$ws.Range("R97").Value = 40.726000728072698
$ws.Range("S97").Value = 1.47343130805293
$ws.Range("T97").Value = 486.05264491258202
$ws.Range("U97").Value = 342.30966638872502
$ws.Range("V97").Value = 322.56594355616602
$ws.Range("W97").Value = 15.2169003382191
$ws.Range("X97").Value = 4.18798463255811
$ws.Range("Y97").Value = 50.405982370350799
$ws.Range("Z97").Value = 23.722500755001999
$ws.Range("AA97").Value = 12.104469221469699

# --- Row 98: "ROW, Fragile States" aggregate recalculation --------------
$ws.Range("C98").Value = 309098
$ws.Range("D98").Value = 463.218566831445
$ws.Range("E98").Value = 12.3684217244649
$ws.Range("F98").Value = 339.28523209115502
$ws.Range("G98").Value = 5.4585414286006797
$ws.Range("H98").Value = 13.1823106611737
$ws.Range("I98").Value = 339093
$ws.Range("J98").Value = 508.169491496467
$ws.Range("K98").Value = 12.613547866275701
$ws.Range("L98").Value = 289.16659480733102
$ws.Range("M98").Value = 4.66566129352994
$ws.Range("N98").Value = 9.7790711468714004
$ws.Range("O98").Value = 43.865928183693299
$ws.Range("P98").Value = 61.477639024278297
$ws.Range("Q98").Value = 129.57615436694499
$ws.Range("R98").Value = 130.90392630404199
$ws.Range("S98").Value = 4.3609665202605798
$ws.Range("T98").Value = 163.844958978216
$ws.Range("U98").Value = 89.073115829803399
$ws.Range("V98").Value = 72.233191160329895
$ws.Range("W98").Value = 46.840676823802298
$ws.Range("X98").Value = 12.6408084530577
$ws.Range("Y98").Value = 155.36580139346199
$ws.Range("Z98").Value = 54.738372205854901
$ws.Range("AA98").Value = 27.358695805456101

# --- Workbook window geometry (best-effort; the runtime's xlsx writer may
#     not persist bookViews state, but attempt it via the COM window model
#     in case it is honoured) ------------------------------------------
try {
    $win = $excel.ActiveWindow
    $win.Left = 57480
    $win.Top = -120
    $win.Width = 19440
    $win.Height = 15000
} catch {
}
